$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain text (preserve exact
# formatting, e.g. trailing zeros / multi-dot "thousands" strings) are forced
# to Text format before the value is written, matching how Excel is driven
# by hand to stop it from auto-converting "1.002"-style strings into numbers.

$ws.Range('D2').Value = '29.163.20'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').Value = '1.851.71'
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.45'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6794'
$ws.Range('E6').Value = '  -5.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07675'
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3032'
$ws.Range('E9').Value = '  -3.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.08'
$ws.Range('E10').Value = '  -6.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08145'
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.842.97'
$ws.Range('E12').Value = '  -3.98%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7208'
$ws.Range('E13').Value = '  -3.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.176'
$ws.Range('E14').Value = '  -2.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.13'
$ws.Range('E15').Value = '  -3.66%  '
$ws.Range('D16').Value = '29.165.45'
$ws.Range('E16').Value = '  -2.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007807'
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.707'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.15'
$ws.Range('E19').Value = '  -2.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '233.22'
$ws.Range('E20').Value = '  -5.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = '2.100.82'
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.422'
$ws.Range('E24').Value = '  -4.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.58'
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.940'
$ws.Range('E26').Value = '  -3.60%  '
$ws.Range('E27').Value = '  -4.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.99'
$ws.Range('E28').Value = '  -3.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.950'
$ws.Range('E29').Value = '  -2.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.389'
$ws.Range('E30').Value = '  -3.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.513'
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.483'
$ws.Range('E32').Value = '  -2.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.002'
$ws.Range('E33').Value = '  -4.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05157'
$ws.Range('E34').Value = '  -4.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.177'
$ws.Range('E35').Value = '  -4.29%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7019'
$ws.Range('E36').Value = '  -5.14%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.020'
$ws.Range('E37').Value = '  +1.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.673'
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01841'
$ws.Range('E39').Value = '  -4.38%  '
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9093'
$ws.Range('E41').Value = '  +2.96%  '
$ws.Range('D42').Value = '1.096.95'
$ws.Range('E42').Value = '  +5.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.974'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4263'
$ws.Range('E44').Value = '  -4.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.93'
$ws.Range('E45').Value = '  -2.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.18'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.762'
$ws.Range('E48').Value = '  -3.10%  '
$ws.Range('D49').Value = '1.995.46'
$ws.Range('E49').Value = '  -3.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.120'
$ws.Range('E50').Value = '  -5.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.896'
$ws.Range('E51').Value = '  -7.88%  '
